# Daily attendance processing - 2025-12-24 11:31:04
# Updates the "Year 5 / B1-7..B1-12 / GENERAL SURGERY" 24/12/2025 sessions
# (rows 25, 44, 63, 172, 191, 210) from "Not Recorded" to "Recorded", refreshes
# the dependent summary statistics, normalises the "Recorded By" text order,
# and narrows column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write literal text into a cell without Excel's automatic
# number/percentage/date re-interpretation, while preserving the cell's
# current style (no new cellXf is left behind once the workbook settles,
# since the scratch "@" style gets reused for every call below).
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# 1. Column I width 14 -> 10 (column H already stores the target raw width of
#    10, so copy it across rather than fight Excel's character-width rounding)
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# ---------------------------------------------------------------------------
# 2. "Recorded By" column: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# ---------------------------------------------------------------------------
$swapRows = @(2,3,20,21,22,24,39,40,41,43,58,59,60,62,77,78,95,96,113,114,131,132,149,150,167,168,169,171,186,187,188,190,205,206,207,209)
foreach ($r in $swapRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}

# ---------------------------------------------------------------------------
# 3. Headline statistics (K/L columns near the top of the sheet)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 78     # Recorded Sessions: 72 -> 78
$ws.Range("L7").Value = 0      # Missing Sessions: 6 -> 0

Set-LiteralText $ws.Range("L9") "35.1%"    # Coverage %: 32.4% -> 35.1%
$ws.Range("L10").Copy() | Out-Null
$ws.Range("L9").PasteSpecial(-4122) | Out-Null

Set-LiteralText $ws.Range("L10") "77.5%"   # Average Attendance %: 77.4% -> 77.5%
$ws.Range("L9").Copy() | Out-Null
$ws.Range("L10").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 4. Per-group statistics block (rows 16-26, columns O/P/R/S) for the groups
#    whose 24/12/2025 session just became Recorded: B1-10, B1-11, B1-12,
#    B1-7, B1-8, B1-9 (summary rows 16,17,18,24,25,26 respectively).
# ---------------------------------------------------------------------------
$statRows = @(
    @{ Row = 16; O = 7; P = 0; R = "36.8%"; S = "73.7%" },
    @{ Row = 17; O = 7; P = 0; R = "36.8%"; S = "58.6%" },
    @{ Row = 18; O = 7; P = 0; R = "36.8%"; S = "83.7%" },
    @{ Row = 24; O = 7; P = 0; R = "36.8%"; S = "69.8%" },
    @{ Row = 25; O = 7; P = 0; R = "36.8%"; S = "73.9%" },
    @{ Row = 26; O = 7; P = 0; R = "36.8%"; S = "69.5%" }
)

foreach ($stat in $statRows) {
    $r = $stat.Row
    $ws.Cells.Item($r, 15).Value = $stat.O   # column O
    $ws.Cells.Item($r, 16).Value = $stat.P   # column P

    $rCell = $ws.Cells.Item($r, 18)          # column R
    Set-LiteralText $rCell $stat.R
    $ws.Range("L10").Copy() | Out-Null
    $rCell.PasteSpecial(-4122) | Out-Null

    $sCell = $ws.Cells.Item($r, 19)          # column S
    Set-LiteralText $sCell $stat.S
    $ws.Range("L10").Copy() | Out-Null
    $sCell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 5. The six session rows that moved from "Not Recorded"/empty to "Recorded".
#    Copy the already-"Recorded" row 6 formatting (green fill) across A:I,
#    then fill in the Recorded-By / Students / Status values.
# ---------------------------------------------------------------------------
$sessionRows = @(
    @{ Row = 25;  G = "dnasr281@gmail.com"; H = "22/31" },
    @{ Row = 44;  G = "dnasr281@gmail.com"; H = "13/19" },
    @{ Row = 63;  G = "dnasr281@gmail.com"; H = "19/21" },
    @{ Row = 172; G = "dnasr281@gmail.com"; H = "24/27" },
    @{ Row = 191; G = "dnasr281@gmail.com"; H = "22/29" },
    @{ Row = 210; G = "dnasr281@gmail.com"; H = "22/29" }
)

$ws.Range("A6:I6").Copy() | Out-Null
foreach ($s in $sessionRows) {
    $r = $s.Row
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122) | Out-Null
}

foreach ($s in $sessionRows) {
    $r = $s.Row
    $ws.Cells.Item($r, 7).Value = $s.G        # G: Recorded By
    $ws.Cells.Item($r, 8).Value = $s.H        # H: Students
    $ws.Cells.Item($r, 9).Value = "Recorded"  # I: Status
}

Write-Host "Daily attendance processing complete"
